{"js": "// Replace the five narrative text segments in the San Diego GRP report body.\n// The paragraph is a single run whose <w:t> text chunks are separated by\n// <w:br/><w:br/> pairs; we locate each old chunk with a body-wide search\n// (exact text, case-sensitive) and replace it in place so the breaks /\n// run formatting around it stay untouched.\nconst replacements = [\n  [\"The economic landscape of San Diego County, as reflected in its Gross Regional Product (GRP), offers a compelling narrative of growth and resilience over recent years. The GRP data, which serves as a vital indicator of economic health, reveals a trajectory of steady expansion from 2019 to 2023. In 2019, San Diego County's GRP stood at approximately $244.28 billion, with a per capita GRP of $73,347. This figure saw a modest increase in 2020, reaching $244.82 billion, despite the challenges posed by the global pandemic. The per capita GRP also rose slightly to $74,278, indicating a resilient economic performance amidst global uncertainties.\", \"The economic landscape of San Diego County has been marked by a dynamic evolution in its Gross Regional Product (GRP) over recent years. This analysis delves into the GRP data from 2019 to 2023, providing insights into the economic vitality of the region, with comparisons to state and national levels.\"],\n  [\"The subsequent years marked a period of robust growth for the county. By 2021, the GRP had surged to $268.87 billion, with a significant rise in per capita GRP to $82,100. This upward trend continued into 2022, with the GRP climbing to $296.68 billion and the per capita figure reaching $90,557. The year 2023 further solidified this growth trajectory, as the GRP peaked at $308.71 billion, and the per capita GRP increased to $94,916. These figures underscore San Diego County's dynamic economic environment, characterized by a diverse range of industries contributing to its overall economic vitality.\", \"In 2019, San Diego County's GRP stood at approximately $244.28 billion, with a per capita GRP of $73,347. This figure slightly increased in 2020 to $244.82 billion, despite a decrease in population, resulting in a per capita GRP of $74,278. The year 2021 marked a significant economic upturn, with the GRP rising to $268.87 billion and the per capita GRP reaching $82,100. This upward trajectory continued into 2022, with the GRP climbing to $296.68 billion and a per capita GRP of $90,557. By 2023, the GRP further increased to $308.71 billion, with a per capita GRP of $94,916, reflecting a robust economic growth pattern.\"],\n  [\"When compared to the state of California and the national economy, San Diego County's performance is noteworthy. California's GRP also experienced growth, rising from approximately $2.99 trillion in 2019 to $3.65 trillion in 2023. The state's per capita GRP followed a similar pattern, increasing from $75,789 in 2019 to $93,800 in 2023. On a national scale, the United States saw its GRP grow from $20.93 trillion in 2019 to $25.96 trillion in 2023, with per capita figures rising from $63,754 to $77,366 over the same period.\", \"When comparing these figures to the state of California, we observe a similar trend. California's GRP in 2019 was approximately $2.99 trillion, with a per capita GRP of $75,789. The state's GRP slightly decreased in 2020 to $2.96 trillion, with a per capita GRP of $74,964. However, like San Diego County, California experienced economic growth in the subsequent years, with the GRP reaching $3.31 trillion in 2021 and $3.54 trillion in 2022, culminating in $3.65 trillion in 2023. The per capita GRP for California followed suit, increasing to $84,587 in 2021, $90,636 in 2022, and $93,800 in 2023.\"],\n  [\"San Diego County's economic growth is driven by a diverse array of industries, each contributing significantly to the region's GRP. Key sectors include government, manufacturing, and professional, scientific, and technical services, which have consistently been major contributors to the county's economic output. The government sector alone accounted for over $52.92 billion in 2023, while manufacturing contributed approximately $31.67 billion. The professional, scientific, and technical services sector also played a crucial role, with a contribution of $37.04 billion in the same year.\", \"On a national scale, the United States' GRP was approximately $20.93 trillion in 2019, with a per capita GRP of $63,754. The national GRP saw a slight decline in 2020 to $20.61 trillion, with a per capita GRP of $62,157. However, the economy rebounded in 2021, with the GRP rising to $22.86 trillion and a per capita GRP of $68,858. This growth continued into 2022 and 2023, with the GRP reaching $24.96 trillion and $25.96 trillion, respectively, and the per capita GRP increasing to $74,889 in 2022 and $77,366 in 2023.\"],\n  [\"In summary, San Diego County's economic performance from 2019 to 2023 reflects a resilient and growing economy, with significant contributions from various sectors. The county's GRP growth, both in absolute terms and on a per capita basis, highlights its economic strength and adaptability in the face of global challenges. As the county continues to build on this foundation, it remains a vital component of California's and the nation's economic landscape.\", \"The data indicates that San Diego County's economic performance has been strong, with consistent growth in GRP and per capita GRP over the years. This growth aligns with the broader economic trends observed at the state and national levels, underscoring the region's resilience and capacity for economic expansion. As we look to the future, San Diego County's economic trajectory suggests a promising outlook, driven by its diverse industries and strategic economic initiatives.\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find expected text: \" + oldText.substring(0, 40));\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the five narrative text segments in the San Diego GRP report body.\n# The paragraph is a single run whose text chunks (<w:t> runs) are separated\n# by line-break pairs; for each chunk we Find the exact old wording (search\n# only, no Replacement text) so the match range collapses onto just that\n# chunk, then assign the new wording straight to Range.Text. Setting the\n# range's Text directly -- rather than using Find's Replace mode -- keeps\n# the straight apostrophes/quotes intact instead of Word's smart-quote\n# autocorrect kicking in, and leaves the surrounding <w:br/> breaks alone.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('The economic landscape of San Diego County, as reflected in its Gross Regional Product (GRP), offers a compelling narrative of growth and resilience over recent years. The GRP data, which serves as a vital indicator of economic health, reveals a trajectory of steady expansion from 2019 to 2023. In 2019, San Diego County''s GRP stood at approximately $244.28 billion, with a per capita GRP of $73,347. This figure saw a modest increase in 2020, reaching $244.82 billion, despite the challenges posed by the global pandemic. The per capita GRP also rose slightly to $74,278, indicating a resilient economic performance amidst global uncertainties.', 'The economic landscape of San Diego County has been marked by a dynamic evolution in its Gross Regional Product (GRP) over recent years. This analysis delves into the GRP data from 2019 to 2023, providing insights into the economic vitality of the region, with comparisons to state and national levels.'),\n    @('The subsequent years marked a period of robust growth for the county. By 2021, the GRP had surged to $268.87 billion, with a significant rise in per capita GRP to $82,100. This upward trend continued into 2022, with the GRP climbing to $296.68 billion and the per capita figure reaching $90,557. The year 2023 further solidified this growth trajectory, as the GRP peaked at $308.71 billion, and the per capita GRP increased to $94,916. These figures underscore San Diego County''s dynamic economic environment, characterized by a diverse range of industries contributing to its overall economic vitality.', 'In 2019, San Diego County''s GRP stood at approximately $244.28 billion, with a per capita GRP of $73,347. This figure slightly increased in 2020 to $244.82 billion, despite a decrease in population, resulting in a per capita GRP of $74,278. The year 2021 marked a significant economic upturn, with the GRP rising to $268.87 billion and the per capita GRP reaching $82,100. This upward trajectory continued into 2022, with the GRP climbing to $296.68 billion and a per capita GRP of $90,557. By 2023, the GRP further increased to $308.71 billion, with a per capita GRP of $94,916, reflecting a robust economic growth pattern.'),\n    @('When compared to the state of California and the national economy, San Diego County''s performance is noteworthy. California''s GRP also experienced growth, rising from approximately $2.99 trillion in 2019 to $3.65 trillion in 2023. The state''s per capita GRP followed a similar pattern, increasing from $75,789 in 2019 to $93,800 in 2023. On a national scale, the United States saw its GRP grow from $20.93 trillion in 2019 to $25.96 trillion in 2023, with per capita figures rising from $63,754 to $77,366 over the same period.', 'When comparing these figures to the state of California, we observe a similar trend. California''s GRP in 2019 was approximately $2.99 trillion, with a per capita GRP of $75,789. The state''s GRP slightly decreased in 2020 to $2.96 trillion, with a per capita GRP of $74,964. However, like San Diego County, California experienced economic growth in the subsequent years, with the GRP reaching $3.31 trillion in 2021 and $3.54 trillion in 2022, culminating in $3.65 trillion in 2023. The per capita GRP for California followed suit, increasing to $84,587 in 2021, $90,636 in 2022, and $93,800 in 2023.'),\n    @('San Diego County''s economic growth is driven by a diverse array of industries, each contributing significantly to the region''s GRP. Key sectors include government, manufacturing, and professional, scientific, and technical services, which have consistently been major contributors to the county''s economic output. The government sector alone accounted for over $52.92 billion in 2023, while manufacturing contributed approximately $31.67 billion. The professional, scientific, and technical services sector also played a crucial role, with a contribution of $37.04 billion in the same year.', 'On a national scale, the United States'' GRP was approximately $20.93 trillion in 2019, with a per capita GRP of $63,754. The national GRP saw a slight decline in 2020 to $20.61 trillion, with a per capita GRP of $62,157. However, the economy rebounded in 2021, with the GRP rising to $22.86 trillion and a per capita GRP of $68,858. This growth continued into 2022 and 2023, with the GRP reaching $24.96 trillion and $25.96 trillion, respectively, and the per capita GRP increasing to $74,889 in 2022 and $77,366 in 2023.'),\n    @('In summary, San Diego County''s economic performance from 2019 to 2023 reflects a resilient and growing economy, with significant contributions from various sectors. The county''s GRP growth, both in absolute terms and on a per capita basis, highlights its economic strength and adaptability in the face of global challenges. As the county continues to build on this foundation, it remains a vital component of California''s and the nation''s economic landscape.', 'The data indicates that San Diego County''s economic performance has been strong, with consistent growth in GRP and per capita GRP over the years. This growth aligns with the broader economic trends observed at the state and national levels, underscoring the region''s resilience and capacity for economic expansion. As we look to the future, San Diego County''s economic trajectory suggests a promising outlook, driven by its diverse industries and strategic economic initiatives.'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Could not find expected text: $($oldText.Substring(0, 40))\"\n    }\n    $range.Text = $newText\n}\n"}
